$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column C with header "Assertion" and first data row "Added to Cart"
$ws.Range("C1").Value = "Assertion"
$ws.Range("C2").Value = "Added to Cart"

# Match the column width used for the new column C (~35.18 chars)
$ws.Range("C1").ColumnWidth = 34.3

# Move selection to C8 to match final saved view state
$ws.Range("C8").Select()
